# Fix formatting issues introduced when scrapping floating point numbers
# (Argentine-style "1.234,56" thousands/decimal separators were mis-scraped;
# normalize to plain "1234.56") and repair two "Razon social" entries where a
# blanket "." -> "" / "," -> "." clean-up pass had mangled the text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H ("Importe") holds these numbers as TEXT (shared strings), not
# real numbers. Assigning a numeric-looking string straight to .Value would
# make Excel coerce it into an actual number (dropping the ".00"), so we
# mark the range as Text first, write the corrected strings, then restore
# the original "Normal" style (General format) so no extra cell formatting
# is introduced.
$importeRange = $ws.Range("H2:H277")
$importeRange.NumberFormat = "@"

$ws.Cells.Item(2, 8).Value = "3600.00"
$ws.Cells.Item(3, 8).Value = "6600.00"
$ws.Cells.Item(4, 8).Value = "19830.00"
$ws.Cells.Item(5, 8).Value = "3350.00"
$ws.Cells.Item(6, 8).Value = "70000.00"
$ws.Cells.Item(7, 8).Value = "320000.00"
$ws.Cells.Item(8, 8).Value = "70000.00"
$ws.Cells.Item(9, 8).Value = "70000.00"
$ws.Cells.Item(10, 8).Value = "105000.00"
$ws.Cells.Item(11, 8).Value = "717000.00"
$ws.Cells.Item(12, 8).Value = "153400.00"
$ws.Cells.Item(13, 8).Value = "1056260.00"
$ws.Cells.Item(14, 8).Value = "1296.00"
$ws.Cells.Item(15, 8).Value = "6050.00"
$ws.Cells.Item(16, 8).Value = "550.00"
$ws.Cells.Item(17, 8).Value = "49998.32"
$ws.Cells.Item(18, 8).Value = "527417.29"
$ws.Cells.Item(19, 8).Value = "148355.08"
$ws.Cells.Item(20, 8).Value = "31271.20"
$ws.Cells.Item(21, 8).Value = "122944.94"
$ws.Cells.Item(22, 8).Value = "10520.00"
$ws.Cells.Item(23, 8).Value = "3453.95"
$ws.Cells.Item(24, 8).Value = "16500.00"
$ws.Cells.Item(25, 8).Value = "68.70"
$ws.Cells.Item(26, 8).Value = "23655.20"
$ws.Cells.Item(27, 8).Value = "5405.00"
$ws.Cells.Item(28, 8).Value = "24069.33"
$ws.Cells.Item(29, 8).Value = "13399.63"
$ws.Cells.Item(30, 8).Value = "2000.00"
$ws.Cells.Item(31, 8).Value = "3250.00"
$ws.Cells.Item(32, 8).Value = "2490.00"
$ws.Cells.Item(33, 8).Value = "1749.02"
$ws.Cells.Item(34, 8).Value = "59900.00"
$ws.Cells.Item(35, 8).Value = "16120.00"
$ws.Cells.Item(36, 8).Value = "300.00"
$ws.Cells.Item(37, 8).Value = "6400.00"
$ws.Cells.Item(38, 8).Value = "90.60"
$ws.Cells.Item(39, 8).Value = "65590.50"
$ws.Cells.Item(40, 8).Value = "1530.00"
$ws.Cells.Item(41, 8).Value = "38902.09"
$ws.Cells.Item(42, 8).Value = "187600.00"
$ws.Cells.Item(43, 8).Value = "496200.20"
$ws.Cells.Item(44, 8).Value = "35352.30"
$ws.Cells.Item(45, 8).Value = "27527.50"
$ws.Cells.Item(46, 8).Value = "46020.00"
$ws.Cells.Item(47, 8).Value = "365503.00"
$ws.Cells.Item(48, 8).Value = "2040.62"
$ws.Cells.Item(49, 8).Value = "31128.50"
$ws.Cells.Item(50, 8).Value = "122.60"
$ws.Cells.Item(51, 8).Value = "3110.74"
$ws.Cells.Item(52, 8).Value = "1517.38"
$ws.Cells.Item(53, 8).Value = "453.75"
$ws.Cells.Item(54, 8).Value = "6010.05"
$ws.Cells.Item(55, 8).Value = "51696.45"
$ws.Cells.Item(56, 8).Value = "1743.00"
$ws.Cells.Item(57, 8).Value = "7473.80"
$ws.Cells.Item(58, 8).Value = "680.00"
$ws.Cells.Item(59, 8).Value = "1108.48"
$ws.Cells.Item(60, 8).Value = "1424.25"
$ws.Cells.Item(61, 8).Value = "2317.00"
$ws.Cells.Item(62, 8).Value = "13857.92"
$ws.Cells.Item(63, 8).Value = "762.30"
$ws.Cells.Item(64, 8).Value = "2220.00"
$ws.Cells.Item(65, 8).Value = "5235.81"
$ws.Cells.Item(66, 8).Value = "20.00"
$ws.Cells.Item(67, 8).Value = "4222.00"
$ws.Cells.Item(68, 8).Value = "9424.00"
$ws.Cells.Item(69, 8).Value = "128.00"
$ws.Cells.Item(70, 8).Value = "20766.03"
$ws.Cells.Item(71, 8).Value = "4126.10"
$ws.Cells.Item(72, 8).Value = "6750.00"
$ws.Cells.Item(73, 8).Value = "64.85"
$ws.Cells.Item(74, 8).Value = "1209.60"
$ws.Cells.Item(75, 8).Value = "37665.54"
$ws.Cells.Item(76, 8).Value = "6405.60"
$ws.Cells.Item(77, 8).Value = "2991.82"
$ws.Cells.Item(78, 8).Value = "3560.00"
$ws.Cells.Item(79, 8).Value = "3750.00"
$ws.Cells.Item(80, 8).Value = "10120.00"
$ws.Cells.Item(81, 8).Value = "8172.00"
$ws.Cells.Item(82, 8).Value = "24591.00"
$ws.Cells.Item(83, 8).Value = "1413.27"
$ws.Cells.Item(84, 8).Value = "757.70"
$ws.Cells.Item(85, 8).Value = "14095.00"
$ws.Cells.Item(86, 8).Value = "5450.00"
$ws.Cells.Item(87, 8).Value = "7505.00"
$ws.Cells.Item(88, 8).Value = "1230.00"
$ws.Cells.Item(89, 8).Value = "2600.00"
$ws.Cells.Item(90, 8).Value = "760.00"
$ws.Cells.Item(91, 8).Value = "10328.00"
$ws.Cells.Item(92, 8).Value = "25200.00"
$ws.Cells.Item(93, 8).Value = "3800.00"
$ws.Cells.Item(94, 8).Value = "11209.00"
$ws.Cells.Item(95, 8).Value = "62260.00"
$ws.Cells.Item(96, 8).Value = "712.05"
$ws.Cells.Item(97, 8).Value = "340.00"
$ws.Cells.Item(98, 8).Value = "840.00"
$ws.Cells.Item(99, 8).Value = "35620.00"
$ws.Cells.Item(100, 8).Value = "3894.80"
$ws.Cells.Item(101, 8).Value = "4810.00"
$ws.Cells.Item(102, 8).Value = "6850.00"
$ws.Cells.Item(103, 8).Value = "526097.59"
$ws.Cells.Item(104, 8).Value = "53261.49"
$ws.Cells.Item(105, 8).Value = "3450.00"
$ws.Cells.Item(106, 8).Value = "15000.00"
$ws.Cells.Item(107, 8).Value = "42.58"
$ws.Cells.Item(108, 8).Value = "365.00"
$ws.Cells.Item(109, 8).Value = "28638.27"
$ws.Cells.Item(110, 8).Value = "124.80"
$ws.Cells.Item(111, 8).Value = "120.00"
$ws.Cells.Item(112, 8).Value = "134.00"
$ws.Cells.Item(113, 8).Value = "7600.00"
$ws.Cells.Item(114, 8).Value = "219000.00"
$ws.Cells.Item(115, 8).Value = "2328.40"
$ws.Cells.Item(116, 8).Value = "689.60"
$ws.Cells.Item(117, 8).Value = "219.12"
$ws.Cells.Item(118, 8).Value = "5200.00"
$ws.Cells.Item(119, 8).Value = "6250.00"
$ws.Cells.Item(120, 8).Value = "100.00"
$ws.Cells.Item(121, 8).Value = "536.00"
$ws.Cells.Item(122, 8).Value = "8660.00"
$ws.Cells.Item(123, 8).Value = "3410.00"
$ws.Cells.Item(124, 8).Value = "29910.11"
$ws.Cells.Item(125, 8).Value = "10839.95"
$ws.Cells.Item(126, 8).Value = "2959.10"
$ws.Cells.Item(127, 8).Value = "1044.78"
$ws.Cells.Item(128, 8).Value = "167.40"
$ws.Cells.Item(129, 8).Value = "3388.00"
$ws.Cells.Item(130, 8).Value = "6472.80"
$ws.Cells.Item(131, 8).Value = "155.00"
$ws.Cells.Item(132, 8).Value = "9600.00"
$ws.Cells.Item(133, 8).Value = "120.00"
$ws.Cells.Item(134, 8).Value = "500.00"
$ws.Cells.Item(135, 8).Value = "4980.00"
$ws.Cells.Item(136, 8).Value = "5242.60"
$ws.Cells.Item(137, 8).Value = "10084.12"
$ws.Cells.Item(138, 8).Value = "8160.00"
$ws.Cells.Item(139, 8).Value = "34100.00"
$ws.Cells.Item(140, 8).Value = "29541.00"
$ws.Cells.Item(141, 8).Value = "21577.00"
$ws.Cells.Item(142, 8).Value = "3898.00"
$ws.Cells.Item(143, 8).Value = "33299.89"
$ws.Cells.Item(144, 8).Value = "5881.00"
$ws.Cells.Item(145, 8).Value = "52000.00"
$ws.Cells.Item(146, 8).Value = "4464.00"
$ws.Cells.Item(147, 8).Value = "20013.00"
$ws.Cells.Item(148, 8).Value = "44800.00"
$ws.Cells.Item(149, 8).Value = "28500.00"
$ws.Cells.Item(150, 8).Value = "1000.00"
$ws.Cells.Item(151, 8).Value = "7200.00"
$ws.Cells.Item(152, 8).Value = "15000.00"
$ws.Cells.Item(153, 8).Value = "6300.00"
$ws.Cells.Item(154, 8).Value = "121250.00"
$ws.Cells.Item(155, 8).Value = "6000.00"
$ws.Cells.Item(156, 8).Value = "8500.00"
$ws.Cells.Item(157, 8).Value = "13789.38"
$ws.Cells.Item(158, 8).Value = "2149.00"
$ws.Cells.Item(159, 8).Value = "1063.20"
$ws.Cells.Item(160, 8).Value = "6790.00"
$ws.Cells.Item(161, 8).Value = "2024.50"
$ws.Cells.Item(162, 8).Value = "2805000.00"
$ws.Cells.Item(163, 8).Value = "15000.00"
$ws.Cells.Item(164, 8).Value = "8000.00"
$ws.Cells.Item(165, 8).Value = "3500.00"
$ws.Cells.Item(166, 8).Value = "4000.00"
$ws.Cells.Item(167, 8).Value = "25978.18"
$ws.Cells.Item(168, 8).Value = "2500.00"
$ws.Cells.Item(169, 8).Value = "3204.50"
$ws.Cells.Item(170, 8).Value = "2556.00"
$ws.Cells.Item(171, 8).Value = "3000.00"
$ws.Cells.Item(172, 8).Value = "3000.00"
$ws.Cells.Item(173, 8).Value = "2000.00"
$ws.Cells.Item(174, 8).Value = "1500.00"
$ws.Cells.Item(175, 8).Value = "7500.00"
$ws.Cells.Item(176, 8).Value = "4000.00"
$ws.Cells.Item(177, 8).Value = "2500.00"
$ws.Cells.Item(178, 8).Value = "4100.00"
$ws.Cells.Item(179, 8).Value = "2500.00"
$ws.Cells.Item(180, 8).Value = "1500.00"
$ws.Cells.Item(181, 8).Value = "17880.00"
$ws.Cells.Item(182, 8).Value = "2000.00"
$ws.Cells.Item(183, 8).Value = "1279.00"
$ws.Cells.Item(184, 8).Value = "250.00"
$ws.Cells.Item(185, 8).Value = "42124.50"
$ws.Cells.Item(186, 8).Value = "6900.00"
$ws.Cells.Item(187, 8).Value = "13200.00"
$ws.Cells.Item(188, 8).Value = "7955.00"
$ws.Cells.Item(189, 8).Value = "1153.32"
$ws.Cells.Item(190, 8).Value = "1400.00"
$ws.Cells.Item(191, 8).Value = "1200.00"
$ws.Cells.Item(192, 8).Value = "15.00"
$ws.Cells.Item(193, 8).Value = "1480.00"
$ws.Cells.Item(194, 8).Value = "9715.00"
$ws.Cells.Item(195, 8).Value = "6350.00"
$ws.Cells.Item(196, 8).Value = "43870.00"
$ws.Cells.Item(197, 8).Value = "228.48"
$ws.Cells.Item(198, 8).Value = "2365.00"
$ws.Cells.Item(199, 8).Value = "2100.00"
$ws.Cells.Item(200, 8).Value = "1827.82"
$ws.Cells.Item(201, 8).Value = "7381.00"
$ws.Cells.Item(202, 8).Value = "79.61"
$ws.Cells.Item(203, 8).Value = "3722.00"
$ws.Cells.Item(204, 8).Value = "40880.00"
$ws.Cells.Item(205, 8).Value = "950.00"
$ws.Cells.Item(206, 8).Value = "1868.78"
$ws.Cells.Item(207, 8).Value = "7200.00"
$ws.Cells.Item(208, 8).Value = "84.00"
$ws.Cells.Item(209, 8).Value = "21735.00"
$ws.Cells.Item(210, 8).Value = "1110.00"
$ws.Cells.Item(211, 8).Value = "21450.00"
$ws.Cells.Item(212, 8).Value = "1590.00"
$ws.Cells.Item(213, 8).Value = "469.24"
$ws.Cells.Item(214, 8).Value = "3300.00"
$ws.Cells.Item(215, 8).Value = "38194.34"
$ws.Cells.Item(216, 8).Value = "25000.00"
$ws.Cells.Item(217, 8).Value = "25000.00"
$ws.Cells.Item(218, 8).Value = "8500.00"
$ws.Cells.Item(219, 8).Value = "25000.00"
$ws.Cells.Item(220, 8).Value = "25000.00"
$ws.Cells.Item(221, 8).Value = "25000.00"
$ws.Cells.Item(222, 8).Value = "54170.00"
$ws.Cells.Item(223, 8).Value = "50000.00"
$ws.Cells.Item(224, 8).Value = "25000.00"
$ws.Cells.Item(225, 8).Value = "42000.00"
$ws.Cells.Item(226, 8).Value = "2550.00"
$ws.Cells.Item(227, 8).Value = "8961.03"
$ws.Cells.Item(228, 8).Value = "1844105.39"
$ws.Cells.Item(229, 8).Value = "3140.00"
$ws.Cells.Item(230, 8).Value = "2350.00"
$ws.Cells.Item(231, 8).Value = "2000.00"
$ws.Cells.Item(232, 8).Value = "51093.43"
$ws.Cells.Item(233, 8).Value = "116190.00"
$ws.Cells.Item(234, 8).Value = "122190.00"
$ws.Cells.Item(235, 8).Value = "116190.00"
$ws.Cells.Item(236, 8).Value = "116190.00"
$ws.Cells.Item(237, 8).Value = "121675.00"
$ws.Cells.Item(238, 8).Value = "106190.00"
$ws.Cells.Item(239, 8).Value = "201846.00"
$ws.Cells.Item(240, 8).Value = "200190.00"
$ws.Cells.Item(241, 8).Value = "294690.00"
$ws.Cells.Item(242, 8).Value = "116190.00"
$ws.Cells.Item(243, 8).Value = "116190.00"
$ws.Cells.Item(244, 8).Value = "116190.00"
$ws.Cells.Item(245, 8).Value = "119340.00"
$ws.Cells.Item(246, 8).Value = "116190.00"
$ws.Cells.Item(247, 8).Value = "200190.00"
$ws.Cells.Item(248, 8).Value = "375690.00"
$ws.Cells.Item(249, 8).Value = "200190.00"
$ws.Cells.Item(250, 8).Value = "117432.00"
$ws.Cells.Item(251, 8).Value = "184190.00"
$ws.Cells.Item(252, 8).Value = "116190.00"
$ws.Cells.Item(253, 8).Value = "116190.00"
$ws.Cells.Item(254, 8).Value = "116190.00"
$ws.Cells.Item(255, 8).Value = "6600000.00"
$ws.Cells.Item(256, 8).Value = "151956.60"
$ws.Cells.Item(257, 8).Value = "226900.00"
$ws.Cells.Item(258, 8).Value = "39700.00"
$ws.Cells.Item(259, 8).Value = "359539.00"
$ws.Cells.Item(260, 8).Value = "60574.00"
$ws.Cells.Item(261, 8).Value = "109024.00"
$ws.Cells.Item(262, 8).Value = "18400.00"
$ws.Cells.Item(263, 8).Value = "2800.00"
$ws.Cells.Item(264, 8).Value = "10000.00"
$ws.Cells.Item(265, 8).Value = "3000.00"
$ws.Cells.Item(266, 8).Value = "213550.00"
$ws.Cells.Item(267, 8).Value = "12200.00"
$ws.Cells.Item(268, 8).Value = "2490.00"
$ws.Cells.Item(269, 8).Value = "1400.00"
$ws.Cells.Item(270, 8).Value = "36023.00"
$ws.Cells.Item(271, 8).Value = "7600.00"
$ws.Cells.Item(272, 8).Value = "99000.00"
$ws.Cells.Item(273, 8).Value = "34500.00"
$ws.Cells.Item(274, 8).Value = "13500.00"
$ws.Cells.Item(275, 8).Value = "3000.00"
$ws.Cells.Item(276, 8).Value = "160.00"
$ws.Cells.Item(277, 8).Value = "1800.00"

$importeRange.Style = "Normal"

# --- Two "Razon social" shared strings were corrupted by the same bad
# find/replace (commas turned into periods, periods stripped) -- restore
# the corrected text exactly as scraped.
$ws.Range("E101").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E187").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E205").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E230").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E268").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
